$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New date column header (column BB, i.e. column 54), matching formatting of column BA (53)
$ws.Cells.Item(1, 54).Value = "24-ago"
$ws.Cells.Item(1, 54).NumberFormat = $ws.Cells.Item(1, 53).NumberFormat

# New data values for the "24-ago" column, one per product row (rows 2-11)
$values = @(15, 11, 10, 14, 10, 13, 12, 18, 25, 14)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 54).Value = $values[$i]
    $ws.Cells.Item($row, 54).HorizontalAlignment = $ws.Cells.Item($row, 53).HorizontalAlignment
    $ws.Cells.Item($row, 54).NumberFormat = $ws.Cells.Item($row, 53).NumberFormat
}

# Move the active selection to the cell below the newly entered column, as in the original edit
[void]$ws.Cells.Item(12, 54).Select()
